$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The document body is a single 20x5 table where every cell holds one
# arithmetic expression ("a+b=" / "a-b="). The diff rewrites the text of
# all 100 cells (each cell keeps its own run formatting). Writing the new
# value straight to each Cell.Range.Text (instead of a document-wide
# Find & Replace) avoids any cross-cell substring collisions, e.g. one
# new value like "25-2=" would otherwise be corrupted by a later
# replacement whose search text ("5-2=") is a substring of it.
$newValues = @(
    @("99-33=", "5+33=", "25-2=", "89-68=", "71-69="),
    @("83-40=", "7+54=", "71-43=", "97-81=", "8+28="),
    @("36+53=", "46-45=", "23+27=", "50-29=", "27+49="),
    @("0+22=", "57-45=", "14+80=", "56+7=", "27-2="),
    @("84-56=", "75-43=", "0+39=", "96-61=", "90-29="),
    @("11+13=", "50+34=", "82-73=", "57+20=", "55+24="),
    @("29+14=", "27-22=", "35+62=", "30-17=", "55+29="),
    @("11+20=", "49+44=", "80+10=", "66-57=", "26+6="),
    @("54+26=", "22+66=", "78-66=", "96-9=", "47-5="),
    @("45-28=", "16+76=", "33+44=", "0+76=", "53-41="),
    @("85-45=", "56+20=", "0+75=", "46-1=", "75+5="),
    @("28+43=", "14+24=", "34-17=", "44+0=", "42+31="),
    @("92-2=", "54+33=", "51-28=", "27+66=", "13+45="),
    @("19-5=", "5+90=", "55-21=", "98-89=", "11-5="),
    @("97-16=", "45+37=", "12+2=", "26+38=", "52-44="),
    @("14+85=", "88-41=", "4-2=", "68-27=", "82-63="),
    @("72-12=", "4+9=", "6+1=", "63+11=", "57+14="),
    @("5-0=", "0+34=", "54+10=", "11+81=", "83-59="),
    @("30-16=", "24+47=", "12+72=", "82+2=", "47+9="),
    @("44+32=", "36-24=", "85-19=", "96-87=", "34-15="),
)

for ($r = 1; $r -le $newValues.Count; $r++) {
    $row = $newValues[$r - 1]
    for ($c = 1; $c -le $row.Count; $c++) {
        $t.Cell($r, $c).Range.Text = $row[$c - 1]
    }
}

Write-Output "Updated $($newValues.Count * $newValues[0].Count) cells"
